$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edits -------------------------------------------------
# Insert 3 new columns before column E (old D shifts to G, etc.)
$ws.Range("E1:G1").EntireColumn.Insert()
# Insert 14 new rows before row 19
$ws.Range("A19:A32").EntireRow.Insert()

# --- Header row (row 1) -------------------------------------------------
$ws.Range("A1").Value = "Test Case"
$ws.Range("B1").Value = "Test x"
$ws.Range("C1").Value = "mu (mean)"
$ws.Range("D1").Value = "sigma (standard deviation)"
$ws.Range("E1").Value = "z-score (expected)"
$ws.Range("F1").Value = "z-score (actual)"
$ws.Range("G1").Value = "Pass?
(Does Expected = Actual?)"

# --- Rows 2-4: summary/test rows ----------------------------------------
$ws.Range("A2").Value = "Column A"
$ws.Range("B2").Value = 2
$ws.Range("C2").Formula = "=AVERAGE(A20:A22)"
$ws.Range("D2").Formula = "=_xlfn.STDEV.P(A20:A22)"
$ws.Range("E2").Formula = "=(B2-C2)/D2"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "Pass"

$ws.Range("A3").Value = "Column B"
$ws.Range("B3").Value = 4
$ws.Range("C3").Formula = "=AVERAGE(B20:B29)"
$ws.Range("D3").Formula = "=_xlfn.STDEV.P(B20:B29)"
$ws.Range("E3").Formula = "=(B3-C3)/D3"
$ws.Range("F3").Value = 0.69631062382279096
$ws.Range("G3").Value = "Pass"

$ws.Range("A4").Value = "Column C"
$ws.Range("B4").Value = 3
$ws.Range("C4").Formula = "=AVERAGE(C20:C32)"
$ws.Range("D4").Formula = "=_xlfn.STDEV.P(C20:C32)"
$ws.Range("E4").Formula = "=(B4-C4)/D4"
$ws.Range("F4").Value = -0.29465944404836197
$ws.Range("G4").Value = "Pass"

# --- Row 12 labels moved from B12:C12 to E12:F12 -------------------------
$ws.Range("B12").ClearContents()
$ws.Range("C12").ClearContents()
$ws.Range("E12").Value = "(results are from spreadsheet formulas)"
$ws.Range("F12").Value = "(results are from the Python function)"

# --- Row 19 column headers for raw data table -----------------------------
$ws.Range("A19").Value = "Column A"
$ws.Range("B19").Value = "Column B"
$ws.Range("C19").Value = "Column C"

# --- Rows 20-32 raw data ---------------------------------------------------
$ws.Range("A20").Value = 1
$ws.Range("B20").Value = -1
$ws.Range("C20").Value = 1

$ws.Range("A21").Value = 2
$ws.Range("B21").Value = -3
$ws.Range("C21").Value = 5

$ws.Range("A22").Value = 3
$ws.Range("B22").Value = -5
$ws.Range("C22").Value = 11

$ws.Range("B23").Value = -7
$ws.Range("C23").Value = 74

$ws.Range("B24").Value = -9
$ws.Range("C24").Value = 62

$ws.Range("B25").Value = 9
$ws.Range("C25").Value = 33

$ws.Range("B26").Value = 7
$ws.Range("C26").Value = 8

$ws.Range("B27").Value = 5
$ws.Range("C27").Value = -99

$ws.Range("B28").Value = 3
$ws.Range("C28").Value = 5

$ws.Range("B29").Value = 1
$ws.Range("C29").Value = 2

$ws.Range("C30").Value = -7
$ws.Range("C31").Value = 64
$ws.Range("C32").Value = 42
